# Commit: "add missing trunc for div in excel."
#
# The sheet computes, every 19 rows, a division `Z(n-1)/AD(n)` that is
# supposed to be an integer "digit shift" in the day-24 ALU simulation, but
# Excel's plain division leaves fractional remainders that then pollute all
# of the downstream running totals. The fix wraps each of those divisions in
# TRUNC(...,0) so the result is floored to an integer, matching the intended
# integer division semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Every Z-cell on the sheet whose formula is a bare "Z../AD.." division -
# replace it with an equivalent TRUNC(...,0) wrapped formula.
$targetRows = @(10, 29, 48, 67, 86, 105, 124, 143, 162, 181, 200, 219, 238, 257)

foreach ($row in $targetRows) {
    $cell = $ws.Range("Z$row")
    $numerator = "Z" + ($row - 1)
    $denominator = "AD$row"
    $cell.Formula = "=TRUNC($numerator/$denominator,0)"
}

# Reflect the author's final cursor position/selection on the sheet.
$ws.Activate()
$ws.Range("Z257").Select()
